$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix capitalization of student name "абдулова таисия" -> "Абдулова таисия"
$ws.Range("B4").Value = "Абдулова таисия"

# 2. Fill in grade values that were added
$ws.Range("F10").Value = 5
$ws.Range("I10").Value = 5

$ws.Range("E24").Value = 5
$ws.Range("F24").Value = 5
$ws.Range("G24").Value = 5
$ws.Range("H24").Value = 5
$ws.Range("I24").Value = 5

# 3. Update the selection/scroll position to match the latest view state
$ws.Cells.Item(10, 9).Select()
